$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="29.469.35"; E="  +0.16%  " },
    @{ r=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.891.06"; E="  -1.48%  " },
    @{ r=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.006"; E="  -0.31%  " },
    @{ r=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="339.83"; E="  +4.76%  " },
    @{ r=6; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.004"; E="  -0.27%  " },
    @{ r=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.4749"; E="  -1.39%  " },
    @{ r=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.3998"; E="  -1.61%  " },
    @{ r=9; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="47.27"; E="  -1.41%  " },
    @{ r=10; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.08025"; E="  -2.45%  " },
    @{ r=11; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.9900"; E="  -1.93%  " },
    @{ r=12; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="23.08"; E="  -1.06%  " },
    @{ r=13; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.869.51"; E="  -1.26%  " },
    @{ r=14; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.936"; E="  -2.40%  " },
    @{ r=15; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.088"; E="  -1.93%  " },
    @{ r=16; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="89.05"; E="  -2.86%  " },
    @{ r=17; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.06774"; E="  -1.12%  " },
    @{ r=18; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.006"; E="  -0.23%  " },
    @{ r=19; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.00001020"; E="  -1.92%  " },
    @{ r=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="17.28"; E="  -1.78%  " },
    @{ r=21; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.004"; E="  -0.32%  " },
    @{ r=22; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="29.465.61"; E="  +0.14%  " },
    @{ r=23; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.496"; E="  -3.01%  " },
    @{ r=24; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="11.64"; E="  -1.45%  " },
    @{ r=25; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.151"; E="  -1.52%  " },
    @{ r=26; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.072.67"; E="  -2.43%  " },
    @{ r=27; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="157.42"; E="  +1.08%  " },
    @{ r=28; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="6.490"; E="  -1.14%  " },
    @{ r=29; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="19.63"; E="  -1.85%  " },
    @{ r=30; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.047"; E="  -3.22%  " },
    @{ r=31; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="118.85"; E="  -1.50%  " },
    @{ r=32; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.9971"; E="  -2.00%  " },
    @{ r=33; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.09524"; E="  -1.30%  " },
    @{ r=34; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.472"; E="  -2.59%  " },
    @{ r=35; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.537"; E="  -0.45%  " },
    @{ r=36; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.382"; E="  +0.51%  " },
    @{ r=37; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.06415"; E="  +5.06%  " },
    @{ r=38; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02240"; E="  -2.07%  " },
    @{ r=39; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.197"; E="  +1.51%  " },
    @{ r=40; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.5821"; E="  -2.27%  " },
    @{ r=41; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="10.52"; E="  -3.05%  " },
    @{ r=42; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="7.713"; E="  -4.22%  " },
    @{ r=43; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1821"; E="  -1.55%  " },
    @{ r=44; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="2.413"; E="  +1.31%  " },
    @{ r=45; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="1.270"; E="  -0.79%  " },
    @{ r=46; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="12.16"; E="  -1.97%  " },
    @{ r=47; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.5499"; E="  -1.61%  " },
    @{ r=48; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.07336"; E="  -3.46%  " },
    @{ r=49; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="1.950"; E="  -0.13%  " },
    @{ r=50; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="116.22"; E="  -1.86%  " },
    @{ r=51; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.379"; E="  -1.93%  " }
)

foreach ($row in $rows) {
    $ws.Range("B" + $row.r).Value = $row.B
    $ws.Range("C" + $row.r).Value = $row.C
    $dval = $row.D
    if ($dval -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range("D" + $row.r).Value = "'" + $dval
    } else {
        $ws.Range("D" + $row.r).Value = $dval
    }
    $ws.Range("E" + $row.r).Value = $row.E
}
